# Generate Report for Handoff
# Updates the localization-status report: files 827253f4-11c5-4be0-abec-9e127905240a
# and 927edbff-a883-4087-ad65-5f4b84f07fa9 moved from "Handed back: in sync with en-US"
# to "Ready for handoff", with refreshed handoff timestamps and an explanatory
# "Error Detail" note on the per-language sheets because the handback file on the
# source repo is behind the latest commit.

$wb = $excel.ActiveWorkbook

$newStatus = "Ready for handoff"
$newOverviewDate = "2016-08-13 22:34:07"
$newZhHandoffDate = "2016-08-13 22:33:56"
$newDeHandoffDate = "2016-08-13 22:34:07"

$errDetail827253 = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/oltest/blob/116024577fc572e872ccc58bcc881c42fa2fe840/e2e/827253f4-11c5-4be0-abec-9e127905240a.md, latest: https://github.com/OpenLocalizationTestOrg/oltest/blob/62e0b5d0494a48059504f255666abe29a59b5746/e2e/827253f4-11c5-4be0-abec-9e127905240a.md."
$errDetail927edbff = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/oltest/blob/116024577fc572e872ccc58bcc881c42fa2fe840/e2e/927edbff-a883-4087-ad65-5f4b84f07fa9.md, latest: https://github.com/OpenLocalizationTestOrg/oltest/blob/62e0b5d0494a48059504f255666abe29a59b5746/e2e/927edbff-a883-4087-ad65-5f4b84f07fa9.md."

# --- Overview sheet: rows 4 (827253f4...) and 5 (927edbff...) ---
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E4").Value = $newStatus
$wsOverview.Range("F4").Value = $newStatus
$wsOverview.Range("G4").Value = $newOverviewDate
$wsOverview.Range("E5").Value = $newStatus
$wsOverview.Range("F5").Value = $newStatus
$wsOverview.Range("G5").Value = $newOverviewDate

# --- zh-cn sheet: rows 4 (827253f4...) and 5 (927edbff...) ---
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C4").Value = $newStatus
$wsZhCn.Range("H4").Value = $newZhHandoffDate
$wsZhCn.Range("P4").Value = $errDetail827253

$wsZhCn.Range("C5").Value = $newStatus
$wsZhCn.Range("H5").Value = $newZhHandoffDate
$wsZhCn.Range("P5").Value = $errDetail927edbff

$wsZhCn.Columns.Item(16).ColumnWidth = 40

# --- de-de sheet: rows 4 (827253f4...) and 5 (927edbff...) ---
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C4").Value = $newStatus
$wsDeDe.Range("H4").Value = $newDeHandoffDate
$wsDeDe.Range("P4").Value = $errDetail827253

$wsDeDe.Range("C5").Value = $newStatus
$wsDeDe.Range("H5").Value = $newDeHandoffDate
$wsDeDe.Range("P5").Value = $errDetail927edbff

$wsDeDe.Columns.Item(16).ColumnWidth = 40
